$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain (non-numeric) string: set directly.
# Cells whose new value LOOKS like a pure number (service IDs) must be
# forced to text (matching the original "t=s" shared-string storage used
# throughout this column) by briefly switching NumberFormat to text,
# then clearing the format again so no stray cell style is left behind.

# --- Plain text / status cells ---
$ws.Range("F4").Value = 'Cannot invoke "org.openqa.selenium.WebElement.getText()" because the return value of "connect_OrderProcessNonSPL.P3P.isElementPresent(String)" is null'
$ws.Range("F11").Value = 'Cannot invoke "org.openqa.selenium.WebElement.isDisplayed()" because "element" is null'
$ws.Range("F12").Value = 'no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: headless chrome=119.0.6045.200)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: ''3.141.59'', revision: ''e82be7d358'', time: ''2018-11-14T08:17:03''
System info: host: ''SIPL92'', ip: ''10.212.130.29'', os.name: ''Windows 10'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''19.0.1''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 119.0.6045.200, chrome: {chromedriverVersion: 119.0.6045.105 (38c72552c5e..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:58698}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 9e620668373f11d91f466650b7e55573
*** Element info: {Using=id, value=lblServiceID}'
$ws.Range("E13").Value = 'PASS'
$ws.Range("F13").Value = 'no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: headless chrome=119.0.6045.200)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: ''3.141.59'', revision: ''e82be7d358'', time: ''2018-11-14T08:17:03''
System info: host: ''SIPL92'', ip: ''10.212.130.29'', os.name: ''Windows 10'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''19.0.1''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 119.0.6045.200, chrome: {chromedriverVersion: 119.0.6045.105 (38c72552c5e..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:58698}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 9e620668373f11d91f466650b7e55573
*** Element info: {Using=id, value=lblServiceID}'
$ws.Range("E14").Value = 'PASS'
$ws.Range("F14").Value = 'no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: headless chrome=119.0.6045.200)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: ''3.141.59'', revision: ''e82be7d358'', time: ''2018-11-14T08:17:03''
System info: host: ''SIPL92'', ip: ''10.212.130.29'', os.name: ''Windows 10'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''19.0.1''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 119.0.6045.200, chrome: {chromedriverVersion: 119.0.6045.105 (38c72552c5e..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:58698}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 9e620668373f11d91f466650b7e55573
*** Element info: {Using=id, value=lblServiceID}'
$ws.Range("E24").Value = 'PASS'
$ws.Range("E25").Value = 'PASS'

# --- Numeric-looking Service ID cells (stored as text) ---
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '10221568'
$ws.Range("C2").ClearFormats()

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '10221569'
$ws.Range("C3").ClearFormats()

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '10221579'
$ws.Range("C4").ClearFormats()

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '10221599'
$ws.Range("C5").ClearFormats()

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = '10221733'
$ws.Range("C11").ClearFormats()

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = '10221741'
$ws.Range("C12").ClearFormats()

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '10221744'
$ws.Range("C13").ClearFormats()

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '10221749'
$ws.Range("C14").ClearFormats()

$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = '136382364'
$ws.Range("C24").ClearFormats()
